$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "username"
$ws.Hyperlinks.Add($ws.Range("B1"), "", "", "", "username")

for ($i = 1; $i -le 10; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $i
}

$ws.Range("B2").Select()
